# destek hesaplaması - destek alan kişi 1 olursa transpose t() değiştirildi.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("destek")
$ws.Activate()

# Ad / isim alanlari (test verisi)
$ws.Range("B2").Value = "test"
$ws.Range("G2").Value = "testtest"

# Tarihler (Dogum_Tarihi, Kaza_Tarihi, EsDT)
$ws.Range("D2").Value = "6/5/1943"
$ws.Range("E2").Value = "3/6/2024"
$ws.Range("H2").Value = "6/2/1964"

# Cocuksay -> 0, cocuk1 bilgileri silindi
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "7/29/2024"

# Anne bilgileri "Yok" oldu
$ws.Range("T2").Value = "Yok"
$ws.Range("U2").Value = "-"
$ws.Range("V2").Value = "7/29/2024"

# Baba bilgileri "Yok" oldu
$ws.Range("W2").Value = "Yok"
$ws.Range("X2").Value = "-"
$ws.Range("Y2").Value = "7/29/2024"

# Eskiden 11. satirda kalmis fazladan veri temizlendi
$ws.Range("W11").ClearContents()

# Gorunum ayarlari
$excel.ActiveWindow.Zoom = 80
$ws.Range("S10").Select()
